# Apply the "Black and White" slot review text edits described in the
# commit.  Word's Find/Replace (and Range.Text assignment) regenerate the
# run list of the paragraph being edited, which silently drops any
# pre-existing empty <w:r/> runs that sit at the start of a paragraph.
# To keep the document structurally identical (aside from the intended
# text edits) we instead use Range.InsertXML with a tightly scoped Range
# (the paragraph's text, excluding its paragraph mark) so only the
# targeted run's text content is replaced.

$d = $word.ActiveDocument

function Set-ParagraphText {
    param($Index, $OldText, $NewText, $RunPropsXml = "")

    $p = $d.Paragraphs.Item($Index)
    $current = $p.Range.Text
    # Paragraph range text includes the trailing paragraph mark (\r) - strip it.
    $trimmed = $current.TrimEnd([char]13)
    if ($trimmed -cne $OldText) {
        throw "Paragraph $Index text mismatch. Expected [$OldText] but found [$trimmed]"
    }

    $start = $p.Range.Start
    $end = $p.Range.End - 1   # exclude the paragraph mark
    $r = $d.Range($start, $end)

    $escaped = $NewText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    $needsPreserve = ($NewText -ne $NewText.Trim())
    if ($needsPreserve) {
        $tOpen = "<w:t xml:space=`"preserve`">"
    } else {
        $tOpen = "<w:t>"
    }

    $runXml = "<w:r>" + $RunPropsXml + $tOpen + $escaped + "</w:t></w:r>"
    $fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($fragment)
}

Set-ParagraphText 1 "Play Black and White Slot Game for Free" `
    "Play Black and White Slot Game Free - Review"

Set-ParagraphText 37 "Unique black and white characters" `
    "Medieval theme with unique black and white characters"

Set-ParagraphText 38 "Free spins feature" `
    "Free spins and expandable grid gameplay features"

Set-ParagraphText 39 "Expandable grid feature" `
    "Charming potential for immersing players in the game"

Set-ParagraphText 40 "High RTP" `
    "High RTP of 98.08%"

Set-ParagraphText 42 "Mediocre game backdrop" `
    "Mediocre game backdrop and unconvincing presence of electric energy"

Set-ParagraphText 43 "Lack of gameplay features" `
    "Limited gameplay features and lack of development in black and white theme"

Set-ParagraphText 44 "Play Black and White Slot Game for Free" `
    "Play Black and White Slot Game Free - Review" `
    "<w:rPr><w:b/></w:rPr>"

Set-ParagraphText 45 `
    "Discover the medieval-themed Black and White slot game with unique black and white characters, free spins, and an expandable grid. Play for free now." `
    "Discover the medieval-themed Black and White slot game with unique characters. Play for free and read our review." `
    "<w:rPr><w:i/></w:rPr>"
